$wb = $excel.ActiveWorkbook

# --- winterToursTest: no longer the active tab; whole-sheet selection -----
$wsWinter = $wb.Worksheets.Item(1)
$wsWinter.Range("A1:XFD1048576").Select() | Out-Null

# --- add the new "summerToursTest" sheet, right after winterToursTest -----
$wsSummer = $wb.Worksheets.Add($null, $wsWinter)
$wsSummer.Name = "summerToursTest"

# Clone winterToursTest's header/row formatting (border, text number format,
# hyperlink style, ...) onto the new sheet before filling in the new values.
$wsWinter.Range("A1:D2").Copy()
$wsSummer.Range("A1:D2").PasteSpecial(-4122)

# Header row (same layout as winterToursTest)
$wsSummer.Range("A1").Value = "Name"
$wsSummer.Range("B1").Value = "Email"
$wsSummer.Range("C1").Value = "Phone"
$wsSummer.Range("D1").Value = "Message"

# Data row for the summer tour contact
$wsSummer.Range("A2").Value = "Amanullah Akbar Ali"
$wsSummer.Range("B2").Value = "amanullah.a@gmail.com"
$wsSummer.Range("C2").Value = "9943357865"
$wsSummer.Range("D2").Value = "Hi This is Amanullah from Erode"

# Hyperlink on the email cell
$wsSummer.Hyperlinks.Add($wsSummer.Range("B2"), "mailto:amanullah.a@gmail.com")

# Re-apply the cloned formatting to B2 so it keeps the Hyperlink cell style
# (with border + text format) instead of the Hyperlinks.Add default style.
$wsWinter.Range("B2").Copy()
$wsSummer.Range("B2").PasteSpecial(-4122)

# Column widths matching winterToursTest's layout
$wsSummer.Columns.Item(1).ColumnWidth = 18.77734375
$wsSummer.Columns.Item(2).ColumnWidth = 23
$wsSummer.Columns.Item(3).ColumnWidth = 11
$wsSummer.Columns.Item(4).ColumnWidth = 30.44140625

# New sheet becomes the active tab, with C12 selected on it
$wsSummer.Range("C12").Select() | Out-Null
$wsSummer.Activate()
